# [MOD] Removed image from the database and changed the app accordingly
#
# This workbook holds a "Concetto / Costrutto / Volume" ER-style reference
# table split across two side-by-side blocks (A:C and E:G). The edit:
#  - renames the "Operatore" concept to "Membro"
#  - removes the "Guida" relationship row entirely (row 17)
#  - rescales a handful of volume estimates
#  - narrows column D and shifts the remembered selection from I9 to J9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Operatore" -> "Membro" (A16)
$ws.Range("A16").Value = "Membro"

# Rescale volume estimates
$ws.Range("C12").Value = 30
$ws.Range("C14").Value = 30
$ws.Range("C15").Value = 300
$ws.Range("C16").Value = 300

# Delete row 17 ("Guida | R | 8000") entirely, shifting rows 18+ up by one
[void]$ws.Rows.Item(17).Delete()

# Column D is narrower now (7.58 -> 5.09)
$ws.Columns.Item(4).ColumnWidth = 5.09

# Remembered selection moves from I9 to J9
[void]$ws.Range("J9").Select()
